# V 0.3 - Momento Fletor em ponto
#
# Updates the structural-analysis workbook so a bending moment can be
# evaluated at an intermediate point of the beam:
#   - AbaProp  : property constants rescaled (100/10000/10 -> 0.01/100000/0.01)
#   - AbaNos   : node list gains a split point at x = 0.01 (was x = 7) and a
#                new end node at x = 10.01
#   - AbaBarras: bar list gains a third bar (node 2 -> node 3) using the same
#                rescaled properties
#   - AbaForca : the point load's moment arm is resolved into a pair of
#                +/-16666.67 couple forces on two rows instead of one
#   - AbaApoio : the second support shifts from node 2 to node 3 and both
#                supports become fully restrained

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AbaProp - property constants
# ---------------------------------------------------------------------------
$wsProp = $wb.Worksheets.Item("AbaProp")
$wsProp.Range("A1").Value = 0.01
$wsProp.Range("B1").Value = 100000
$wsProp.Range("C1").Value = 0.01

# ---------------------------------------------------------------------------
# AbaNos - nodes
# ---------------------------------------------------------------------------
$wsNos = $wb.Worksheets.Item("AbaNos")
$wsNos.Range("A2").Value = 0.01
$wsNos.Range("B2").Value = 0

$wsNos.Range("A4").Value = 10.01
$wsNos.Range("B4").Value = 0

# ---------------------------------------------------------------------------
# AbaBarras - bars
# ---------------------------------------------------------------------------
$wsBarras = $wb.Worksheets.Item("AbaBarras")
$wsBarras.Range("C1").Value = 0.01
$wsBarras.Range("D1").Value = 100000
$wsBarras.Range("E1").Value = 0.01

$wsBarras.Range("C2").Value = 0.01
$wsBarras.Range("D2").Value = 100000
$wsBarras.Range("E2").Value = 0.01

$wsBarras.Range("A3").Value = 2
$wsBarras.Range("B3").Value = 3
$wsBarras.Range("C3").Value = 0.01
$wsBarras.Range("D3").Value = 100000
$wsBarras.Range("E3").Value = 0.01

# ---------------------------------------------------------------------------
# AbaForca - forces
# ---------------------------------------------------------------------------
$wsForca = $wb.Worksheets.Item("AbaForca")
$wsForca.Range("D1").Value = -16666.67

$wsForca.Range("A2").Value = 2
$wsForca.Range("B2").Value = 0
$wsForca.Range("C2").Value = -10000
$wsForca.Range("D2").Value = 16666.67

# ---------------------------------------------------------------------------
# AbaApoio - supports
# ---------------------------------------------------------------------------
$wsApoio = $wb.Worksheets.Item("AbaApoio")
$wsApoio.Range("B1").Value = $true
$wsApoio.Range("D1").Value = $true

$wsApoio.Range("A2").Value = 3
$wsApoio.Range("D2").Value = $true
